# 10Th - MB for single stock and added new group
# Insert 3 new date columns (C:E) ahead of the existing data block (old C:F -> F:I),
# stamp the new "Jun_27" date in column B, carry the previous date ("Jun_26")
# into the 3 newly inserted header cells, and backfill the new data columns
# with the "UN" placeholder rating used throughout the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the label that used to sit in B1 (the most-recent date column)
# before we overwrite it with the new date.
$previousDateLabel = $ws.Range("B1").Value2

# Insert three blank columns at C:E - this shifts the old C:F block to F:I.
$ws.Range("C1:E1").EntireColumn.Insert()

# New header row: B1 becomes the newest date; the freshly inserted C1:E1
# headers pick up the date that used to be in B1.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1:E1").Value = $previousDateLabel

# Backfill the new C:E columns on every used data row with the "UN" placeholder.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("C" + $r + ":E" + $r).Value = "UN"
}
